$wb = $excel.ActiveWorkbook

# Put some content on Sheet2 ("Not empty"), then select A2 and activate the sheet
# so it becomes the active tab shown when the workbook is opened.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "Not empty"

$ws2.Activate()
$ws2.Range("A2").Select()
